$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a "DD-MM-YYYY" style date string into a cell as literal text,
# without letting Excel auto-convert it into a date serial number (which it
# does for ambiguous values such as "01-09-2025" where the day could also be
# read as a month). Temporarily switching the cell to Text format forces the
# literal string to be kept, then the format is switched back to General so
# the cell matches the rest of the (unformatted) data column.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = "General"
}

# --- 1. Grow the table by one row at the bottom (row 8), copying formatting
#        (values + styles) from the current last data row (row 7). ---
$ws.Range("A7:F7").Copy($ws.Range("A8:F8"))

# --- 2. Shift the existing data rows 2-7 down to rows 3-8 one at a time,
#        working from the bottom up so each row is read before it gets
#        overwritten by the row above it. ---
for ($r = 7; $r -ge 2; $r--) {
    $dst = $r + 1
    $ws.Cells.Item($dst, 1).Value = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($dst, 2).Value = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($dst, 3).Value = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($dst, 4).Value = $ws.Cells.Item($r, 4).Value()
    Set-TextValue $ws.Cells.Item($dst, 5) $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($dst, 6).Value = $ws.Cells.Item($r, 6).Value()
}

# --- 3. Write the new price entry into row 2. ---
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "ALUMINIUM INGOT"
$ws.Range("C2").Value = "IE07"
$ws.Range("D2").Value = 275.25
Set-TextValue $ws.Range("E2") "20-09-2025"
$ws.Range("F2").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf"

# --- 4. Rebuild the hyperlinks for column F (rows 2-8) so that each one
#        points at the URL now shown in that row's F cell, and restore the
#        plain (non-hyperlink) look of the cell afterwards. ---
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le 8; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $target = $cell.Value()
    $ws.Hyperlinks.Add($cell, $target) | Out-Null
    $cell.Font.Underline = -4142
    $cell.Font.ThemeColor = 1
}
